$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1881188118811881
$ws.Range("C2").Value = 0.5544554455445545
$ws.Range("J2").Value = 0.0297029702970297
$ws.Range("P2").Value = 0.1353135313531353
$ws.Range("S2").Value = 0.0924092409240924
$ws.Range("B3").Value = 0.01129943502824859
$ws.Range("C3").Value = 0.02824858757062147
$ws.Range("J3").Value = 0.03389830508474576
$ws.Range("P3").Value = 0.7570621468926554
$ws.Range("S3").Value = 0.1694915254237288
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.6190476190476191
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.0873015873015873
$ws.Range("D6").Value = 0.003968253968253968
$ws.Range("E6").Value = 0.003968253968253968
$ws.Range("F6").Value = 0.09523809523809523
$ws.Range("J6").Value = 0.2063492063492063
$ws.Range("O6").Value = 0.01984126984126984
$ws.Range("Q6").Value = 0.1944444444444444
$ws.Range("R6").Value = 0.05158730158730158
$ws.Range("S6").Value = 0.3373015873015873
$ws.Range("B7").Value = 0.1256281407035176
$ws.Range("D7").Value = 0.02010050251256281
$ws.Range("F7").Value = 0.05527638190954774
$ws.Range("J7").Value = 0.07537688442211055
$ws.Range("O7").Value = 0.005025125628140704
$ws.Range("Q7").Value = 0.2160804020100502
$ws.Range("R7").Value = 0.07035175879396985
$ws.Range("S7").Value = 0.4321608040201005
$ws.Range("B8").Value = 0.081374321880651
$ws.Range("D8").Value = 0.02350813743218806
$ws.Range("F8").Value = 0.07414104882459313
$ws.Range("J8").Value = 0.09222423146473779
$ws.Range("O8").Value = 0.009041591320072333
$ws.Range("Q8").Value = 0.2296564195298373
$ws.Range("R8").Value = 0.1030741410488246
$ws.Range("S8").Value = 0.3869801084990959
$ws.Range("B9").Value = 0.08530805687203792
$ws.Range("D9").Value = 0.01895734597156398
$ws.Range("F9").Value = 0.08530805687203792
$ws.Range("J9").Value = 0.07582938388625593
$ws.Range("O9").Value = 0.01421800947867299
$ws.Range("Q9").Value = 0.1611374407582938
$ws.Range("R9").Value = 0.1090047393364929
$ws.Range("S9").Value = 0.4502369668246445
$ws.Range("B10").Value = 0.09799382716049383
$ws.Range("D10").Value = 0.0154320987654321
$ws.Range("E10").Value = 0.0007716049382716049
$ws.Range("F10").Value = 0.05478395061728395
$ws.Range("J10").Value = 0.1126543209876543
$ws.Range("O10").Value = 0.0162037037037037
$ws.Range("Q10").Value = 0.2237654320987654
$ws.Range("R10").Value = 0.1003086419753086
$ws.Range("S10").Value = 0.3780864197530864
$ws.Range("G11").Value = 0.1538461538461539
$ws.Range("J11").Value = 0.0695970695970696
$ws.Range("K11").Value = 0.1941391941391941
$ws.Range("L11").Value = 0.5677655677655677
$ws.Range("S11").Value = 0.01465201465201465
$ws.Range("G12").Value = 0.802547770700637
$ws.Range("J12").Value = 0.1401273885350318
$ws.Range("K12").Value = 0.01910828025477707
$ws.Range("L12").Value = 0.01273885350318471
$ws.Range("S12").Value = 0.02547770700636943
$ws.Range("F13").Value = 0.01612903225806452
$ws.Range("G13").Value = 0.6129032258064516
$ws.Range("J13").Value = 0.3225806451612903
$ws.Range("S13").Value = 0.04838709677419355
$ws.Range("F15").Value = 0.01932367149758454
$ws.Range("H15").Value = 0.2222222222222222
$ws.Range("I15").Value = 0.05314009661835749
$ws.Range("J15").Value = 0.3091787439613526
$ws.Range("K15").Value = 0.04347826086956522
$ws.Range("M15").Value = 0.02415458937198068
$ws.Range("O15").Value = 0.0821256038647343
$ws.Range("S15").Value = 0.2463768115942029
$ws.Range("F16").Value = 0.02061855670103093
$ws.Range("H16").Value = 0.2010309278350516
$ws.Range("I16").Value = 0.09793814432989691
$ws.Range("J16").Value = 0.3814432989690721
$ws.Range("K16").Value = 0.09793814432989691
$ws.Range("M16").Value = 0.04123711340206185
$ws.Range("N16").Value = 0.005154639175257732
$ws.Range("O16").Value = 0.06185567010309279
$ws.Range("S16").Value = 0.09278350515463918
$ws.Range("F17").Value = 0.01462522851919561
$ws.Range("H17").Value = 0.2193784277879342
$ws.Range("I17").Value = 0.08775137111517367
$ws.Range("J17").Value = 0.3985374771480805
$ws.Range("K17").Value = 0.09506398537477148
$ws.Range("M17").Value = 0.01645338208409506
$ws.Range("N17").Value = 0.001828153564899452
$ws.Range("O17").Value = 0.05850091407678245
$ws.Range("S17").Value = 0.1078610603290676
$ws.Range("F18").Value = 0.03347280334728033
$ws.Range("H18").Value = 0.1924686192468619
$ws.Range("I18").Value = 0.08786610878661087
$ws.Range("J18").Value = 0.4267782426778243
$ws.Range("K18").Value = 0.100418410041841
$ws.Range("M18").Value = 0.02092050209205021
$ws.Range("O18").Value = 0.06276150627615062
$ws.Range("S18").Value = 0.07531380753138076
$ws.Range("F19").Value = 0.02339622641509434
$ws.Range("H19").Value = 0.2286792452830189
$ws.Range("I19").Value = 0.08679245283018867
$ws.Range("J19").Value = 0.3743396226415094
$ws.Range("K19").Value = 0.08528301886792453
$ws.Range("M19").Value = 0.03094339622641509
$ws.Range("N19").Value = 0.001509433962264151
$ws.Range("O19").Value = 0.05056603773584906
$ws.Range("S19").Value = 0.1184905660377358
